# Generate Report for Handoff
# Updates the localization-status report with a fresh handoff run:
# new source GUID/file name, new content hash, and new handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "5618eef6-2572-4309-abff-b8645fe5ce31"
$newHash = "eedaa5a4e5c001da6d6e901a393d35f43e221077"

$newFileName = "$newGuid.md"
$newPathName = "e2e\$newGuid.md"

$newOverviewDateTime = "2016-08-15 16:54:18"
$newZhCnHandoffDateTime = "2016-08-15 16:54:13"
$newDeDeHandoffDateTime = "2016-08-15 16:54:18"

$newZhCnXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeDeXlf = "$newGuid.$newHash.de-de.xlf"

# The hyperlink target URL (points at the source .md file in the repo) is
# unchanged by this edit -- only the visible display text / cell text move
# to the new file name.
$fileHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52203eea0d4e48c4ea0d4b3517eea449a156215a/e2e/6c8c8a14-df86-427d-9abb-4e40b101daaa.md"

function Update-Hyperlink($ws, $cellRef, $newDisplay) {
    $cell = $ws.Range($cellRef)
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($cell, $fileHyperlinkAddress, "", "", $newDisplay) | Out-Null
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
Update-Hyperlink $wsOverview "B2" $newPathName
$wsOverview.Range("G2").Value = $newOverviewDateTime

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newFileName
Update-Hyperlink $wsZhCn "A2" $newFileName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDateTime

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newFileName
Update-Hyperlink $wsDeDe "A2" $newFileName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newDeDeHandoffDateTime
